$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the absolute path recorded by Excel for this workbook's folder
$wb.Path = "C:\Users\ckely\Desktop\Inventario\UGBplacas\"

# Update the numeric values for the 280(+) (O) and 200(+) (R) columns
$ws.Range("O2").Value = 640
$ws.Range("O3").Value = 640
$ws.Range("O4").Value = 800
$ws.Range("O5").Value = 504
$ws.Range("O6").Value = 548
$ws.Range("O7").Value = 560
$ws.Range("R2").Value = 120

# Move/change the active selection to R2
$ws.Range("R2").Select()
